# Implement security vulnerability checks
# Adds a new row (row 73) of sensor data to each of the four worksheets,
# mirroring the existing row layout (time, lengths, checksum, decoded values).

$wb = $excel.ActiveWorkbook

function Add-SensorRow {
    param($ws, $row, $time, $col2, $col3, $col4, $col5, $col6, $col7, $col8, $col9)

    $ws.Cells.Item($row, 1).Value = $time
    $ws.Cells.Item($row, 2).Value = $col2
    $ws.Cells.Item($row, 3).Value = $col3
    $ws.Cells.Item($row, 4).Value = $col4
    $ws.Cells.Item($row, 5).Value = $col5
    $ws.Cells.Item($row, 6).Value = $col6

    # Column 7 holds a numeric-looking string that exceeds normal numeric
    # precision, so it must be forced to text to avoid being coerced into a
    # floating point number (losing exact digits). Resetting the style back
    # to Normal keeps the cell looking like its untouched neighbors.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = $col7
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = $col8
    $ws.Cells.Item($row, 9).Value = $col9
}

# Sheet 1: ROW35-FE-LIFTER
$ws1 = $wb.Worksheets.Item(1)
Add-SensorRow $ws1 73 "2025-03-07 08:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

# Sheet 2: ROW35-MID-LIFTER
$ws2 = $wb.Worksheets.Item(2)
Add-SensorRow $ws2 73 "2025-03-07 08:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

# Sheet 3: ROW02-FE-LIFTER
$ws3 = $wb.Worksheets.Item(3)
Add-SensorRow $ws3 73 "2025-03-07 08:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

# Sheet 4: ROW02-MID-LIFTER
$ws4 = $wb.Worksheets.Item(4)
Add-SensorRow $ws4 73 "2025-03-07 08:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3
